# Updates numeric "F"/"G" column figures (e.g. ticket counts / prices) across
# the "展览" (Exhibitions), "演出" (Performances) and "全部类型" (All types)
# sheets, matching the refreshed gh-pages data snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# sheet1 (展览)
$ws1.Range("F2").Value = 1115
$ws1.Range("F3").Value = 240
$ws1.Range("F5").Value = 1764
$ws1.Range("F6").Value = 655
$ws1.Range("F8").Value = 453
$ws1.Range("F9").Value = 4343
$ws1.Range("F10").Value = 47
$ws1.Range("F11").Value = 447
$ws1.Range("F14").Value = 1275
$ws1.Range("F17").Value = 2954
$ws1.Range("F18").Value = 1784
$ws1.Range("F21").Value = 162
$ws1.Range("F23").Value = 915
$ws1.Range("F24").Value = 291
$ws1.Range("F26").Value = 2261
$ws1.Range("F27").Value = 981
$ws1.Range("F28").Value = 2334
$ws1.Range("F29").Value = 239
$ws1.Range("F30").Value = 686
$ws1.Range("F31").Value = 517
$ws1.Range("F33").Value = 876
$ws1.Range("F35").Value = 1074
$ws1.Range("F36").Value = 882
$ws1.Range("F37").Value = 1150
$ws1.Range("F39").Value = 385
$ws1.Range("F40").Value = 506
$ws1.Range("F42").Value = 271
$ws1.Range("F43").Value = 3457

# sheet2 (演出)
$ws2.Range("F22").Value = 31
$ws2.Range("G25").Value = 100

# sheet4 (全部类型)
$ws4.Range("F2").Value = 1115
$ws4.Range("F3").Value = 240
$ws4.Range("F6").Value = 1764
$ws4.Range("F7").Value = 655
$ws4.Range("F9").Value = 453
$ws4.Range("F10").Value = 4343
$ws4.Range("F11").Value = 47
$ws4.Range("F16").Value = 1275
$ws4.Range("F17").Value = 2954
$ws4.Range("F19").Value = 1784
$ws4.Range("F27").Value = 915
$ws4.Range("F28").Value = 291
$ws4.Range("F29").Value = 2261
$ws4.Range("F32").Value = 981
$ws4.Range("F33").Value = 2334
$ws4.Range("F34").Value = 687
$ws4.Range("F35").Value = 517
$ws4.Range("F36").Value = 876
$ws4.Range("F37").Value = 1074
$ws4.Range("F38").Value = 882
$ws4.Range("F39").Value = 1150
$ws4.Range("F40").Value = 385
$ws4.Range("F41").Value = 506
$ws4.Range("F44").Value = 31
$ws4.Range("F47").Value = 271
$ws4.Range("F48").Value = 3457
